$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the new rows first (values are set later, in a
#     specific order, so that the shared-string table is built up in the
#     same sequence the original author typed things in) ---

# Row 21 (2018-02-12): copy the taller "B6/C13/D7" style pattern used by row 9
$ws.Range("B9:D9").Copy()
$ws.Range("B21:D21").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(21).RowHeight = 31.5

# Row 22 (2018-02-12): copy the regular "B6/C14/D7" style pattern used by row 11
$ws.Range("B11:D11").Copy()
$ws.Range("B22:D22").PasteSpecial(-4122)  # xlPasteFormats

# Row 23 (2018-02-13): same regular style pattern
$ws.Range("B11:D11").Copy()
$ws.Range("B23:D23").PasteSpecial(-4122)  # xlPasteFormats

# Dates for the new rows
$ws.Range("B21").Value = 43143
$ws.Range("B22").Value = 43143
$ws.Range("B23").Value = 43144

# --- Fill in the cell text content ---
$ws.Range("D21").Value = "45 min"
$ws.Range("C22").Value = "Création script base de données."
$ws.Range("C21").Value = "Documentation, création de la solution."
$ws.Range("C20").Value = "Rédaction Introduction/analyse."
$ws.Range("C19").Value = "Création du MCD - MLD."
$ws.Range("C23").Value = "Création UI acceuil"
$ws.Range("D23").Value = "1h 30 min"

$ws.Range("D22").Value = "45 min"

# --- Update active selection to D23 ---
[void]$ws.Range("D23").Select()
